# Fix Training Data Issue
# The "Date" column (BF) held a malformed literal like "6-16-2013-14".
# Correct it to the real ISO date text "2014-06-16" for every data row
# (rows 2-31). The cells must stay plain text (not get reinterpreted as a
# real Excel date serial), so we force a Text number format on the range
# before writing the corrected strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Range("BF$row").Value = "2014-06-16"
}
